$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 36 (row 48): DC vs RR
$ws.Range("E48").Value = 30
$ws.Range("H48").Value = 40
$ws.Range("K48").Value = 100
$ws.Range("N48").Value = 50
$ws.Range("Q48").Value = 20
$ws.Range("T48").Value = 0
$ws.Range("W48").Value = 80
$ws.Range("Z48").Value = 60
$ws.Range("AC48").Value = 70

# Contest 37 (row 49): SRH vs PBKS
$ws.Range("E49").Value = 60
$ws.Range("H49").Value = 30
$ws.Range("K49").Value = 80
$ws.Range("N49").Value = 0
$ws.Range("Q49").Value = 20
$ws.Range("T49").Value = 50
$ws.Range("W49").Value = 70
$ws.Range("Z49").Value = 40
$ws.Range("AC49").Value = 100
